$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 95,3
$arr[0,0] = 24.9816047538945
$arr[0,1] = 5.628583713734685
$arr[0,2] = 0
$arr[1,0] = 48.02857225639664
$arr[1,1] = 17.72820822527561
$arr[1,2] = 1
$arr[2,0] = 39.2797576724562
$arr[2,1] = 11.28711962152653
$arr[2,2] = 0
$arr[3,0] = 33.94633936788146
$arr[3,1] = 15.17141382329406
$arr[3,2] = 0
$arr[4,0] = 16.23978081344811
$arr[4,1] = 9.985844582977499
$arr[4,2] = 0
$arr[5,0] = 12.32334448672798
$arr[5,1] = 13.20765846071259
$arr[5,2] = 0
$arr[6,0] = 44.64704583099741
$arr[6,1] = 20.11102277086097
$arr[6,2] = 1
$arr[7,0] = 34.04460046972835
$arr[7,1] = 9.575963309832449
$arr[7,2] = 0
$arr[8,0] = 38.32290311184182
$arr[8,1] = 6.539598196575859
$arr[8,2] = 0
$arr[9,0] = 10.8233797718321
$arr[9,1] = 10.79502905827536
$arr[9,2] = 0
$arr[10,0] = 48.79639408647977
$arr[10,1] = 8.224425745080088
$arr[10,2] = 0
$arr[11,0] = 43.29770563201687
$arr[11,1] = 23.59395304685146
$arr[11,2] = 1
$arr[12,0] = 18.49356442713105
$arr[12,1] = 21.16240759128834
$arr[12,2] = 0
$arr[13,0] = 17.27299868828403
$arr[13,1] = 17.66807513020847
$arr[13,2] = 0
$arr[14,0] = 17.33618039413735
$arr[14,1] = 22.42921180375436
$arr[14,2] = 0
$arr[15,0] = 22.16968971838151
$arr[15,1] = 21.07344153798229
$arr[15,2] = 0
$arr[16,0] = 30.99025726528951
$arr[16,1] = 8.731401177720716
$arr[16,2] = 0
$arr[17,0] = 27.27780074568463
$arr[17,1] = 22.85117996979956
$arr[17,2] = 1
$arr[18,0] = 21.64916560792168
$arr[18,1] = 15.78684483831301
$arr[18,2] = 0
$arr[19,0] = 34.47411578889518
$arr[19,1] = 21.14880310328125
$arr[19,2] = 1
$arr[20,0] = 15.57975442608167
$arr[20,1] = 22.92182599846986
$arr[20,2] = 0
$arr[21,0] = 21.68578594140872
$arr[21,1] = 11.36006949943728
$arr[21,2] = 0
$arr[22,0] = 24.65447373174767
$arr[22,1] = 7.201038490553535
$arr[22,2] = 0
$arr[23,0] = 28.24279936868144
$arr[23,1] = 9.558703250838834
$arr[23,2] = 0
$arr[24,0] = 17.98695128633439
$arr[24,1] = 21.36029531844986
$arr[24,2] = 0
$arr[25,0] = 30.56937753654447
$arr[25,1] = 22.21461166512687
$arr[25,2] = 1
$arr[26,0] = 33.6965827544817
$arr[26,1] = 5.139042610623815
$arr[26,2] = 0
$arr[27,0] = 11.85801650879991
$arr[27,1] = 15.21494605155132
$arr[27,2] = 0
$arr[28,0] = 34.30179407605753
$arr[28,1] = 13.34822006297558
$arr[28,2] = 0
$arr[29,0] = 16.82096494749166
$arr[29,1] = 9.442156209414605
$arr[29,2] = 0
$arr[30,0] = 12.60206371941118
$arr[30,1] = 7.397307346673657
$arr[30,2] = 0
$arr[31,0] = 48.62528132298237
$arr[31,1] = 23.85819407825038
$arr[31,2] = 1
$arr[32,0] = 42.33589392465844
$arr[32,1] = 11.4640586404151
$arr[32,2] = 0
$arr[33,0] = 22.18455076693483
$arr[33,1] = 15.37581243486732
$arr[33,2] = 0
$arr[34,0] = 13.90688456025535
$arr[34,1] = 19.06037917790356
$arr[34,2] = 0
$arr[35,0] = 37.36932106048627
$arr[35,1] = 12.27259204758588
$arr[35,2] = 0
$arr[36,0] = 27.60609974958405
$arr[36,1] = 24.43564165441921
$arr[36,2] = 1
$arr[37,0] = 14.88152939379115
$arr[37,1] = 24.24894589884222
$arr[37,2] = 0
$arr[38,0] = 29.80707640445081
$arr[38,1] = 10.03564591650728
$arr[38,2] = 0
$arr[39,0] = 11.37554084460874
$arr[39,1] = 14.94497011784771
$arr[39,2] = 0
$arr[40,0] = 46.37281608315129
$arr[40,1] = 11.01756619633539
$arr[40,2] = 0
$arr[41,0] = 20.35119926400068
$arr[41,1] = 10.69680988754935
$arr[41,2] = 0
$arr[42,0] = 36.50089137415928
$arr[42,1] = 5.737738947090656
$arr[42,2] = 0
$arr[43,0] = 22.46844304357644
$arr[43,1] = 17.19128667959794
$arr[43,2] = 0
$arr[44,0] = 30.80272084711243
$arr[44,1] = 15.05358046457723
$arr[44,2] = 0
$arr[45,0] = 31.86841117373119
$arr[45,1] = 6.029575024999787
$arr[45,2] = 0
$arr[46,0] = 17.39417822102108
$arr[46,1] = 10.57292928473223
$arr[46,2] = 0
$arr[47,0] = 48.78338511058234
$arr[47,1] = 23.16531771933307
$arr[47,2] = 1
$arr[48,0] = 41.00531293444458
$arr[48,1] = 9.791237813339448
$arr[48,2] = 0
$arr[49,0] = 47.57995766256757
$arr[49,1] = 7.897897441824462
$arr[49,2] = 0
$arr[50,0] = 45.79309401710595
$arr[50,1] = 14.78905520555126
$arr[50,2] = 1
$arr[51,0] = 33.9159991524434
$arr[51,1] = 24.71300908221201
$arr[51,2] = 1
$arr[52,0] = 46.87496940092467
$arr[52,1] = 9.841105430230009
$arr[52,2] = 0
$arr[53,0] = 13.53970008207678
$arr[53,1] = 18.44271094811757
$arr[53,2] = 0
$arr[54,0] = 17.83931449676581
$arr[54,1] = 20.23239230657435
$arr[54,2] = 0
$arr[55,0] = 11.80909155642152
$arr[55,1] = 9.752750879847994
$arr[55,2] = 0
$arr[56,0] = 23.01321323053057
$arr[56,1] = 19.56432697223719
$arr[56,2] = 0
$arr[57,0] = 25.54709158757928
$arr[57,1] = 12.35566265438506
$arr[57,2] = 0
$arr[58,0] = 20.85396127095584
$arr[58,1] = 17.64611661187159
$arr[58,2] = 0
$arr[59,0] = 43.14950036607718
$arr[59,1] = 17.6705942152179
$arr[59,2] = 1
$arr[60,0] = 24.27013306774357
$arr[60,1] = 15.71549368149517
$arr[60,2] = 0
$arr[61,0] = 21.23738038749523
$arr[61,1] = 6.805795401088166
$arr[61,2] = 0
$arr[62,0] = 31.70784332632994
$arr[62,1] = 21.70604991178476
$arr[62,2] = 1
$arr[63,0] = 15.63696899899051
$arr[63,1] = 11.41560129943472
$arr[63,2] = 0
$arr[64,0] = 42.08787923016159
$arr[64,1] = 8.730370207997085
$arr[64,2] = 0
$arr[65,0] = 12.98202574719083
$arr[65,1] = 5.815502831095278
$arr[65,2] = 0
$arr[66,0] = 49.47547746402069
$arr[66,1] = 16.81785886376484
$arr[66,2] = 1
$arr[67,0] = 40.88979077186629
$arr[67,1] = 18.55128723684565
$arr[67,2] = 1
$arr[68,0] = 17.94862726136689
$arr[68,1] = 5.331756578557123
$arr[68,2] = 0
$arr[69,0] = 10.2208846849441
$arr[69,1] = 15.24186116598562
$arr[69,2] = 0
$arr[70,0] = 42.61845713819336
$arr[70,1] = 9.529915503958758
$arr[70,2] = 0
$arr[71,0] = 38.27429375390469
$arr[71,1] = 17.903455808189
$arr[71,2] = 1
$arr[72,0] = 39.16028672163949
$arr[72,1] = 8.487328580099829
$arr[72,2] = 0
$arr[73,0] = 40.85081386743783
$arr[73,1] = 18.81875476204932
$arr[73,2] = 1
$arr[74,0] = 12.96178606936362
$arr[74,1] = 12.73470692601075
$arr[74,2] = 0
$arr[75,0] = 24.33862914177091
$arr[75,1] = 23.73459977473469
$arr[75,2] = 1
$arr[76,0] = 14.63476238100519
$arr[76,1] = 7.750418882919865
$arr[76,2] = 0
$arr[77,0] = 44.52413703502374
$arr[77,1] = 11.82132702100517
$arr[77,2] = 0
$arr[78,0] = 34.93192507310232
$arr[78,1] = 7.269470424811781
$arr[78,2] = 0
$arr[79,0] = 23.23592099410597
$arr[79,1] = 23.49387236557126
$arr[79,2] = 1
$arr[80,0] = 12.54233401144095
$arr[80,1] = 22.54678706761962
$arr[80,2] = 0
$arr[81,0] = 22.43929286862649
$arr[81,1] = 10.15883255430311
$arr[81,2] = 0
$arr[82,0] = 23.00733288106988
$arr[82,1] = 18.19968092068358
$arr[82,2] = 0
$arr[83,0] = 39.18424713352256
$arr[83,1] = 21.34444400402432
$arr[83,2] = 1
$arr[84,0] = 35.50229885420853
$arr[84,1] = 16.10401623198925
$arr[84,2] = 0
$arr[85,0] = 45.48850970305306
$arr[85,1] = 15.59301156712013
$arr[85,2] = 1
$arr[86,0] = 28.88859700647797
$arr[86,1] = 9.837045818009035
$arr[86,2] = 0
$arr[87,0] = 14.78376983753207
$arr[87,1] = 6.862055356117985
$arr[87,2] = 0
$arr[88,0] = 38.5297914889198
$arr[88,1] = 22.94431515906653
$arr[88,2] = 1
$arr[89,0] = 40.4314019446759
$arr[89,1] = 23.00836114326661
$arr[89,2] = 1
$arr[90,0] = 40.83868719818244
$arr[90,1] = 11.78059582097401
$arr[90,2] = 0
$arr[91,0] = 29.75182385457563
$arr[91,1] = 11.98419149225322
$arr[91,2] = 0
$arr[92,0] = 27.10164073434198
$arr[92,1] = 22.94220519905154
$arr[92,2] = 1
$arr[93,0] = 11.01676506976381
$arr[93,1] = 22.74172848530235
$arr[93,2] = 0
$arr[94,0] = 14.31565707973218
$arr[94,1] = 20.59751091715248
$arr[94,2] = 0

$range = $ws.Range("A70:C164")
$range.Value = $arr

Write-Host "Added $($arr.GetLength(0)) rows"
